$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 759135.8
$ws.Range("J112").Value = 808075.2
$ws.Range("L112").Value = 2424225.6
$ws.Range("N112").Value = -2426441.6
$ws.Range("H114").Value = 38500
$ws.Range("J114").Value = 38500
$ws.Range("L114").Value = 38500
$ws.Range("N114").Value = -47178
$ws.Range("H137").Value = 5779.9414
$ws.Range("I137").Value = 4388.7026
$ws.Range("J137").Value = 9456.786
$ws.Range("K137").Value = 13166.1078
$ws.Range("L137").Value = 28370.358
$ws.Range("M137").Value = -10616.1078
$ws.Range("N137").Value = -33470.358
$ws.Range("H138").Value = 1860.3732
$ws.Range("I138").Value = 1512.5358
$ws.Range("J138").Value = 2110.1025
$ws.Range("K138").Value = 4537.607400000001
$ws.Range("L138").Value = 6330.3075
$ws.Range("M138").Value = 602.3925999999992
$ws.Range("N138").Value = -16610.3075
$ws.Range("H139").Value = 61092.8
$ws.Range("J139").Value = 61092.8
$ws.Range("L139").Value = 61092.8
$ws.Range("N139").Value = -71372.8
$ws.Range("H140").Value = 38507.418
$ws.Range("J140").Value = 38507.418
$ws.Range("L140").Value = 38507.418
$ws.Range("N140").Value = -48867.418
$ws.Range("H141").Value = 2482.4614
$ws.Range("I141").Value = 1134.9445
$ws.Range("K141").Value = 3404.8335
$ws.Range("M141").Value = 1775.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9410.493
$ws.Range("I32").Value = 8444.725
$ws.Range("K32").Value = 8444.725
$ws.Range("M32").Value = -8157.725
$ws.Range("H61").Value = 1951.6774
$ws.Range("I61").Value = 1036.85
$ws.Range("J61").Value = 3615
$ws.Range("K61").Value = 1036.85
$ws.Range("L61").Value = 3615
$ws.Range("M61").Value = -824.8499999999999
$ws.Range("N61").Value = -4039
$ws.Range("H123").Value = 32551.2
$ws.Range("J123").Value = 32551.2
$ws.Range("L123").Value = 32551.2
$ws.Range("N123").Value = -42351.2
$ws.Range("H132").Value = 21742174
$ws.Range("I132").Value = 35717270
$ws.Range("J132").Value = 3131.5557
$ws.Range("K132").Value = 107151810
$ws.Range("L132").Value = 9394.667099999999
$ws.Range("M132").Value = -107149280
$ws.Range("N132").Value = -14454.6671
$ws.Range("H133").Value = 36615.25
$ws.Range("J133").Value = 36615.25
$ws.Range("L133").Value = 36615.25
$ws.Range("N133").Value = -41675.25
$ws.Range("H136").Value = 1951.6774
$ws.Range("I136").Value = 1036.85
$ws.Range("J136").Value = 3615
$ws.Range("K136").Value = 3110.55
$ws.Range("L136").Value = 10845
$ws.Range("M136").Value = -560.5499999999997
$ws.Range("N136").Value = -15945

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4169.1665
$ws.Range("I105").Value = 4000.8
$ws.Range("K105").Value = 4000.8
$ws.Range("M105").Value = -2253.8
$ws.Range("H134").Value = 2230.9
$ws.Range("I134").Value = 1811.9474
$ws.Range("J134").Value = 3557.5833
$ws.Range("K134").Value = 5435.8422
$ws.Range("L134").Value = 10672.7499
$ws.Range("M134").Value = -2900.8422
$ws.Range("N134").Value = -15742.7499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3925265.5
$ws.Range("I31").Value = 2036.7241
$ws.Range("J31").Value = 5956937.5
$ws.Range("K31").Value = 2036.7241
$ws.Range("L31").Value = 5956937.5
$ws.Range("M31").Value = -1741.7241
$ws.Range("N31").Value = -5957527.5
$ws.Range("H34").Value = 3925265.5
$ws.Range("I34").Value = 2036.7241
$ws.Range("J34").Value = 5956937.5
$ws.Range("K34").Value = 2036.7241
$ws.Range("L34").Value = 5956937.5
$ws.Range("M34").Value = -1834.7241
$ws.Range("N34").Value = -5957341.5
$ws.Range("H99").Value = 1831.3334
$ws.Range("I99").Value = 1782
$ws.Range("J99").Value = 1885.6
$ws.Range("K99").Value = 1782
$ws.Range("L99").Value = 1885.6
$ws.Range("M99").Value = -284
$ws.Range("N99").Value = -4881.6
$ws.Range("H126").Value = 1831.3334
$ws.Range("I126").Value = 1782
$ws.Range("J126").Value = 1885.6
$ws.Range("K126").Value = 5346
$ws.Range("L126").Value = 5656.799999999999
$ws.Range("M126").Value = -2876
$ws.Range("N126").Value = -10596.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 287.60526
$ws.Range("J12").Value = 405.4
$ws.Range("L12").Value = 1216.2
$ws.Range("N12").Value = -1562.2
$ws.Range("H68").Value = 1320.3478
$ws.Range("I68").Value = 861.4286
$ws.Range("J68").Value = 1705.84
$ws.Range("K68").Value = 2584.2858
$ws.Range("L68").Value = 5117.52
$ws.Range("M68").Value = -1773.2858
$ws.Range("N68").Value = -6739.52
$ws.Range("H71").Value = 1320.3478
$ws.Range("I71").Value = 861.4286
$ws.Range("J71").Value = 1705.84
$ws.Range("K71").Value = 7752.8574
$ws.Range("L71").Value = 15352.56
$ws.Range("M71").Value = -3696.8574
$ws.Range("N71").Value = -23464.56
$ws.Range("H123").Value = 1986.9231
$ws.Range("I123").Value = 1919.1666
$ws.Range("J123").Value = 2800
$ws.Range("K123").Value = 5757.4998
$ws.Range("L123").Value = 8400
$ws.Range("M123").Value = -3307.4998
$ws.Range("N123").Value = -13300
$ws.Range("H131").Value = 48378.277
$ws.Range("I131").Value = 150476.12
$ws.Range("J131").Value = 27435.129
$ws.Range("K131").Value = 451428.36
$ws.Range("L131").Value = 82305.387
$ws.Range("M131").Value = -446388.36
$ws.Range("N131").Value = -92385.387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1294.9231
$ws.Range("I102").Value = 1148
$ws.Range("J102").Value = 1360.2222
$ws.Range("K102").Value = 1148
$ws.Range("L102").Value = 1360.2222
$ws.Range("M102").Value = 474
$ws.Range("N102").Value = -4604.2222
$ws.Range("H126").Value = 12995.5
$ws.Range("I126").Value = 13660.111
$ws.Range("J126").Value = 7014
$ws.Range("K126").Value = 40980.333
$ws.Range("L126").Value = 21042
$ws.Range("M126").Value = -38510.333
$ws.Range("N126").Value = -25982
$ws.Range("H132").Value = 2829.0715
$ws.Range("I132").Value = 1978
$ws.Range("J132").Value = 3566.6667
$ws.Range("K132").Value = 5934
$ws.Range("L132").Value = 10700.0001
$ws.Range("M132").Value = -3404
$ws.Range("N132").Value = -15760.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58826136
$ws.Range("I7").Value = 100001150
$ws.Range("J7").Value = 4686.4287
$ws.Range("K7").Value = 100001150
$ws.Range("L7").Value = 4686.4287
$ws.Range("M7").Value = -100001038
$ws.Range("N7").Value = -4910.4287
$ws.Range("H126").Value = 58826136
$ws.Range("I126").Value = 100001150
$ws.Range("J126").Value = 4686.4287
$ws.Range("K126").Value = 300003450
$ws.Range("L126").Value = 14059.2861
$ws.Range("M126").Value = -300000980
$ws.Range("N126").Value = -18999.2861
$ws.Range("H132").Value = 4387.4
$ws.Range("I132").Value = 4126.4
$ws.Range("J132").Value = 4778.9
$ws.Range("K132").Value = 12379.2
$ws.Range("L132").Value = 14336.7
$ws.Range("M132").Value = -9849.199999999999
$ws.Range("N132").Value = -19396.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 6667676.5
$ws.Range("I107").Value = 942
$ws.Range("J107").Value = 20001146
$ws.Range("K107").Value = 2826
$ws.Range("L107").Value = 60003438
$ws.Range("M107").Value = -906
$ws.Range("N107").Value = -60007278
$ws.Range("H123").Value = 35416.668
$ws.Range("J123").Value = 35416.668
$ws.Range("L123").Value = 35416.668
$ws.Range("N123").Value = -45216.668
$ws.Range("H132").Value = 1922.234
$ws.Range("I132").Value = 1568
$ws.Range("K132").Value = 4704
$ws.Range("M132").Value = -2174
$ws.Range("H136").Value = 271345.44
$ws.Range("I136").Value = 371345.47
$ws.Range("J136").Value = 1345.3
$ws.Range("K136").Value = 1114036.41
$ws.Range("L136").Value = 4035.9
$ws.Range("M136").Value = -1111486.41
$ws.Range("N136").Value = -9135.9
